$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" '58.207.40'
Set-TextValue "E2" '  +2.05%  '

Set-TextValue "D3" '2.352.29'
Set-TextValue "E3" '  +0.15%  '

Set-TextValue "E4" '  -0.06%  '

Set-TextValue "D5" '540.14'
Set-TextValue "E5" '  +1.56%  '

Set-TextValue "D6" '135.95'
Set-TextValue "E6" '  +2.41%  '

Set-TextValue "E7" '  +0.37%  '

Set-TextValue "E8" '  +5.48%  '

Set-TextValue "E9" '  +0.30%  '

Set-TextValue "E10" '  +4.76%  '

Set-TextValue "E11" '  -0.85%  '

Set-TextValue "D12" '0.353'
Set-TextValue "E12" '  +2.08%  '

Set-TextValue "B13" 'Avalanche'
Set-TextValue "C13" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D13" '23.81'
Set-TextValue "E13" '  +1.36%  '

Set-TextValue "B14" 'WrappedliquidstakedEther2.0'
Set-TextValue "C14" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D14" '2.770.40'
Set-TextValue "E14" '  +1.05%  '

Set-TextValue "D15" '58.171.27'
Set-TextValue "E15" '  +1.82%  '

Set-TextValue "D16" '0.0000134'
Set-TextValue "E16" '  +0.54%  '

Set-TextValue "D17" '2.352.21'
Set-TextValue "E17" '  +0.62%  '

Set-TextValue "E18" '  +2.43%  '

Set-TextValue "D19" '331.99'
Set-TextValue "E19" '  -1.65%  '

Set-TextValue "E20" '  +2.49%  '

Set-TextValue "E21" '  -1.30%  '

Set-TextValue "E22" '  -0.01%  '

Set-TextValue "D23" '62.79'
Set-TextValue "E23" '  +1.69%  '

Set-TextValue "E24" '  +0.43%  '

Set-TextValue "B25" 'Binance-PegBSC-USD'
Set-TextValue "C25" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D25" '1.01'
Set-TextValue "E25" '  +1.27%  '

Set-TextValue "B26" 'InternetComputer(DFINITY)'
Set-TextValue "C26" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D26" '8.49'
Set-TextValue "E26" '  -2.73%  '

Set-TextValue "E27" '  +2.93%  '

Set-TextValue "E28" '  +1.54%  '

Set-TextValue "D29" '171.23'
Set-TextValue "E29" '  -1.72%  '

Set-TextValue "E30" '  +1.52%  '

Set-TextValue "E31" '  +0.10%  '

Set-TextValue "E32" '  +12.61%  '

Set-TextValue "D33" '18.43'
Set-TextValue "E33" '  -0.59%  '

Set-TextValue "E34" '  +0.03%  '

Set-TextValue "E35" '  +6.72%  '

Set-TextValue "D36" '0.999'
Set-TextValue "E36" '  +0.53%  '

Set-TextValue "E37" '  -0.55%  '

Set-TextValue "D38" '1.65'
Set-TextValue "E38" '  +4.64%  '

Set-TextValue "D39" '39.22'
Set-TextValue "E39" '  -0.12%  '

Set-TextValue "D40" '144.90'
Set-TextValue "E40" '  -2.98%  '

Set-TextValue "D41" '296.22'
Set-TextValue "E41" '  +4.87%  '

Set-TextValue "E42" '  +0.81%  '

Set-TextValue "E43" '  +1.01%  '

Set-TextValue "D44" '0.0948'
Set-TextValue "E44" '  +1.85%  '

Set-TextValue "E45" '  +1.62%  '

Set-TextValue "E46" '  +0.13%  '

Set-TextValue "E47" '  +0.57%  '

Set-TextValue "E48" '  +1.24%  '

Set-TextValue "E49" '  +0.25%  '

Set-TextValue "E50" '  -0.08%  '

